$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B21").Value = 6220
$ws.Range("C21").Value = 980
$ws.Range("D21").Value = 5600875
$ws.Range("E21").Value = 900.4622186495177
$ws.Range("F21").Value = 7.967366776601281
$ws.Range("G21").Value = 3.375527426160341
$ws.Range("H21").Value = 27.84548725716576
